# Auto update ESO load data
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Updated load values (columns D:AA, i.e. col indices 4..27) for rows 10-14,
# plus the brand-new row 15 (including its A/B/C key columns).
$rowData = @{
    10 = @(4867,4668,4578,4573,4555,4709,5165,5650,6056,6211,6145,6057,5958,5885,5779,5767,5778,5913,6202,6184,5938,5613,5361,4947)
    11 = @(4595,4361,4266,4234,4225,4259,4417,4713,5238,5558,5639,5588,5523,5468,5383,5341,5400,5570,5740,5668,5376,5006,4855,4575)
    12 = @(4320,4108,4011,3960,3985,4018,4158,4462,4888,5254,5413,5435,5414,5347,5263,5246,5349,5531,5742,5675,5418,5083,4828,4439)
    13 = @(4641,4416,4340,4311,4344,4506,4950,5454,5919,6165,6180,6159,6019,5983,5830,5830,5785,5914,6233,6230,5996,5660,5344,4892)
    14 = @(4908,4713,4625,4620,4602,4769,5263,5787,6226,6395,6326,6234,6131,6056,5944,5932,5946,6103,6440,6421,6161,5817,5550,5112)
    15 = @(4980,4786,4697,4693,4675,4842,5336,5860,6300,6468,6396,6300,6191,6112,5995,5982,5996,6160,6511,6492,6232,5889,5623,5187)
}

foreach ($r in $rowData.Keys) {
    $values = $rowData[$r]
    for ($i = 0; $i -lt $values.Length; $i++) {
        $col = 4 + $i   # D is column 4
        $ws.Cells.Item($r, $col).Value = $values[$i]
    }
}

# New row 15 key columns: Year, Month, Day
$ws.Cells.Item(15, 1).Value = 2026
$ws.Cells.Item(15, 2).Value = 2
$ws.Cells.Item(15, 3).Value = 18
